# Delete the 'current_time' row from the 'constants' sheet.
# When Excel re-saves the workbook, the now-unused 'current_time' shared
# string is dropped from the string table and every subsequent shared
# string index (and row index) is compacted/shifted accordingly - exactly
# what the target diff shows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constants")

# Find the row whose column A holds the literal text "current_time".
$found = $ws.Columns.Item(1).Find("current_time", [Type]::Missing, [Type]::Missing, 1)

if ($found -ne $null) {
    $found.EntireRow.Delete()
}
